{"js": "// Update the date stamp and every two-digit multiplication prompt in the\n// worksheet to the new day's values. Every \"before\" string is unique in the\n// document, so a case-sensitive body.search() for each literal string\n// followed by a Range.insertText(..., \"Replace\") on the hit keeps each run's\n// character formatting (font/size) intact while only swapping the text.\n\nconst replacements = [\n  [\"2023-11-08 Wednesday\", \"2023-11-09 Thursday\"],\n  [\"79\u00d759=\", \"62\u00d711=\"],\n  [\"25\u00d791=\", \"32\u00d735=\"],\n  [\"88\u00d793=\", \"38\u00d780=\"],\n  [\"49\u00d764=\", \"43\u00d749=\"],\n  [\"86\u00d715=\", \"29\u00d739=\"],\n  [\"26\u00d758=\", \"18\u00d745=\"],\n  [\"30\u00d728=\", \"14\u00d798=\"],\n  [\"95\u00d763=\", \"33\u00d729=\"],\n  [\"40\u00d724=\", \"83\u00d740=\"],\n  [\"45\u00d798=\", \"29\u00d788=\"],\n  [\"57\u00d732=\", \"95\u00d797=\"],\n  [\"79\u00d734=\", \"73\u00d753=\"],\n  [\"33\u00d779=\", \"75\u00d735=\"],\n  [\"52\u00d725=\", \"61\u00d740=\"],\n  [\"30\u00d759=\", \"51\u00d768=\"],\n  [\"59\u00d774=\", \"96\u00d769=\"],\n  [\"65\u00d776=\", \"80\u00d775=\"],\n  [\"83\u00d764=\", \"34\u00d767=\"],\n  [\"67\u00d715=\", \"46\u00d727=\"],\n  [\"55\u00d728=\", \"91\u00d791=\"],\n  [\"36\u00d787=\", \"82\u00d792=\"],\n  [\"93\u00d738=\", \"26\u00d726=\"],\n  [\"13\u00d716=\", \"72\u00d756=\"],\n  [\"20\u00d739=\", \"67\u00d737=\"],\n  [\"72\u00d718=\", \"47\u00d764=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < found.items.length; i++) {\n    found.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date stamp and every two-digit multiplication prompt in the\n# worksheet to the new day's values. Every \"before\" string is unique in the\n# document, so a plain literal Find/Replace (no wildcards) for each pair is\n# sufficient and keeps each run's character formatting (font/size) intact.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2023-11-08 Wednesday\", \"2023-11-09 Thursday\"),\n    @(\"79\u00d759=\", \"62\u00d711=\"),\n    @(\"25\u00d791=\", \"32\u00d735=\"),\n    @(\"88\u00d793=\", \"38\u00d780=\"),\n    @(\"49\u00d764=\", \"43\u00d749=\"),\n    @(\"86\u00d715=\", \"29\u00d739=\"),\n    @(\"26\u00d758=\", \"18\u00d745=\"),\n    @(\"30\u00d728=\", \"14\u00d798=\"),\n    @(\"95\u00d763=\", \"33\u00d729=\"),\n    @(\"40\u00d724=\", \"83\u00d740=\"),\n    @(\"45\u00d798=\", \"29\u00d788=\"),\n    @(\"57\u00d732=\", \"95\u00d797=\"),\n    @(\"79\u00d734=\", \"73\u00d753=\"),\n    @(\"33\u00d779=\", \"75\u00d735=\"),\n    @(\"52\u00d725=\", \"61\u00d740=\"),\n    @(\"30\u00d759=\", \"51\u00d768=\"),\n    @(\"59\u00d774=\", \"96\u00d769=\"),\n    @(\"65\u00d776=\", \"80\u00d775=\"),\n    @(\"83\u00d764=\", \"34\u00d767=\"),\n    @(\"67\u00d715=\", \"46\u00d727=\"),\n    @(\"55\u00d728=\", \"91\u00d791=\"),\n    @(\"36\u00d787=\", \"82\u00d792=\"),\n    @(\"93\u00d738=\", \"26\u00d726=\"),\n    @(\"13\u00d716=\", \"72\u00d756=\"),\n    @(\"20\u00d739=\", \"67\u00d737=\"),\n    @(\"72\u00d718=\", \"47\u00d764=\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $new\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
